$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15
$ws.Range("A15").Value = 112388168
$ws.Range("B15").Value = 96735
$ws.Range("D15").Value = 'VU'
$ws.Range("E15").Value = 220787
$ws.Range("F15").Value = 'Knärot'
$ws.Range("G15").Value = 'Goodyera repens'
$ws.Range("H15").Value = '(L.) R. Br.'
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = '50'
$ws.Range("J15").NumberFormat = "@"
$ws.Range("J15").Value = 'plantor/tuvor'
$ws.Range("P15").Value = 'Övre Tjärna, Dlr'
$ws.Range("Q15").Value = 520808
$ws.Range("R15").Value = 6706139
$ws.Range("Z15").Value = ""
$ws.Range("AB15").Value = ""
$ws.Range("AW15").Value = 'Holger Martinussen'
$ws.Range("AX15").Value = 'Holger Martinussen, Uno Skog, Anton Björk'
# Row 16
$ws.Range("A16").Value = 112389296
$ws.Range("B16").Value = 94301
$ws.Range("D16").Value = 'NT'
$ws.Range("E16").Value = 53
$ws.Range("F16").Value = 'Vedtrappmossa'
$ws.Range("G16").Value = 'Crossocalyx hellerianus'
$ws.Range("H16").Value = '(Nees ex Lindenb.) Meyl.'
$ws.Range("I16").Value = ""
$ws.Range("J16").Value = ""
$ws.Range("P16").Value = 'Paradiset, Dlr'
$ws.Range("Q16").Value = 520702
$ws.Range("R16").Value = 6706232
$ws.Range("Z16").NumberFormat = "@"
$ws.Range("Z16").Value = '10:03'
$ws.Range("AB16").NumberFormat = "@"
$ws.Range("AB16").Value = '10:03'
$ws.Range("AW16").Value = 'Uno Skog'
$ws.Range("AX16").Value = 'Uno Skog, Holger Martinussen, Anton Björk'
# Row 17
$ws.Range("B17").Value = 96735
# Row 18
$ws.Range("A18").Value = 112435720
$ws.Range("B18").Value = 83506
$ws.Range("E18").Value = 241
$ws.Range("F18").Value = 'Gransotdyna'
$ws.Range("G18").Value = 'Camarops tubulina'
$ws.Range("H18").Value = '(Alb. & Schwein.:Fr.) Shear'
$ws.Range("Q18").Value = 520679
$ws.Range("R18").Value = 6706241
$ws.Range("S18").Value = 10
# Row 19
$ws.Range("A19").Value = 112435727
$ws.Range("B19").Value = 89553
$ws.Range("E19").Value = 1202
$ws.Range("F19").Value = 'Ullticka'
$ws.Range("G19").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H19").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q19").Value = 520887
$ws.Range("R19").Value = 6706018
$ws.Range("S19").Value = 5
# Row 20
$ws.Range("A20").Value = 112389127
$ws.Range("B20").Value = 96735
$ws.Range("Q20").Value = 520688
$ws.Range("R20").Value = 6706226
$ws.Range("Z20").NumberFormat = "@"
$ws.Range("Z20").Value = '09:47'
$ws.Range("AB20").NumberFormat = "@"
$ws.Range("AB20").Value = '09:47'
# Row 21
$ws.Range("A21").Value = 112388247
$ws.Range("B21").Value = 96735
$ws.Range("Q21").Value = 520819
$ws.Range("R21").Value = 6706134
$ws.Range("Z21").NumberFormat = "@"
$ws.Range("Z21").Value = '08:46'
$ws.Range("AB21").NumberFormat = "@"
$ws.Range("AB21").Value = '08:46'
# Row 22
$ws.Range("A22").Value = 112388456
$ws.Range("B22").Value = 96735
$ws.Range("Q22").Value = 520759
$ws.Range("R22").Value = 6706151
$ws.Range("Z22").NumberFormat = "@"
$ws.Range("Z22").Value = '09:00'
$ws.Range("AB22").NumberFormat = "@"
$ws.Range("AB22").Value = '09:00'
# Row 23
$ws.Range("B23").Value = 90835
